# compare haplotypes with filtered ED50s
#
# Insert a new "mtORF" column before the existing CBASS column (this shifts
# CBASS/TPC/lat/lon one column to the right, e.g. I->J, J->K, K->L, L->M),
# then populate the mtORF column with haplotype calls, and populate a new
# "TPC2" column (which lands in the previously-empty column N) with the
# filtered ED50 flag.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift CBASS/TPC/lat/lon right by inserting the new mtORF column at I.
$ws.Columns.Item(9).Insert()

# Headers
$ws.Cells.Item(1, 9).Value = "mtORF"
$ws.Cells.Item(1, 14).Value = "TPC2"

# mtORF (column I) haplotype calls, by row
$mtorf = @{
  9  = "3b"
  10 = "1a"
  11 = "1a"
  12 = "1a"
  13 = "1a"
  15 = "3b"
  16 = "1a"
  17 = "1a"
  18 = "3b"
  19 = "1a"
  20 = "1a"
  21 = "3b"
  22 = "3b"
  23 = "3b"
  24 = "1a"
  25 = "3b"
  26 = "1a"
  27 = "1a"
}

foreach ($r in $mtorf.Keys) {
  $ws.Cells.Item($r, 9).Value = $mtorf[$r]
}

# TPC2 (column N) filtered-ED50 flags, by row
$tpc2 = @{
  9  = "Y"
  10 = "Y"
  11 = "Y"
  12 = "Y"
  13 = "Y"
  15 = "Y"
  16 = "Y"
  17 = "Y"
  18 = "Y"
  20 = "Y"
  21 = "Y"
  22 = "Y"
  23 = "Y"
  24 = "Y"
  25 = "Y"
}

foreach ($r in $tpc2.Keys) {
  $ws.Cells.Item($r, 14).Value = $tpc2[$r]
}

# Match the author's final selection (cosmetic, but reflected in the diff).
$ws.Range("N26").Select() | Out-Null
